$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the filter description text (new UO 25917 added to the list)
$ws.Range("A4").Value2 = "({Item Informação} = 56:PAGAMENTOS TOTAIS (EXERCICIO E RAP)) E ({Órgão UGE - Orçam. Fiscal S/N} = PERTENCE) E ({Unidade Orçamentária} = 40904:FUNDO DO REGIME GERAL DA PREVID.SOCIAL- FRGPS, 55902:FUNDO DO REGIME GERAL DA PREVID.SOCIAL-FRGPS, 33904:FUNDO DO REGIME GERAL DA PREVIDENCIA SOCIAL, 25917:FUNDO DO REGIME GERAL DE PREVIDENCIA SOCIAL) E ({Ano Lançamento} ({Número Ano}) = 2019)"

# Update the "Movimento R$" values in column U for the replaced data rows
$ws.Range("U12").Value2 = 2194900.8
$ws.Range("U13").Value2 = 6161800449.26
$ws.Range("U14").Value2 = 57482832.8299999
$ws.Range("U15").Value2 = 58038967.9499998
$ws.Range("U16").Value2 = 9251715.32999992

$ws.Range("U18").Value2 = 572934223.48
$ws.Range("U19").Value2 = 956666330.04
$ws.Range("U20").Value2 = 876046353.92
$ws.Range("U21").Value2 = 940334079.67
$ws.Range("U22").Value2 = 25452262.8099999

$ws.Range("U25").Value2 = 483053256.24
$ws.Range("U26").Value2 = 255969667.98
$ws.Range("U27").Value2 = 242092703.06
$ws.Range("U28").Value2 = 283156554.82

$ws.Range("U30").Value2 = 45604902847.01
$ws.Range("U31").Value2 = 27345470660.24
$ws.Range("U32").Value2 = 36735225735.85
$ws.Range("U33").Value2 = 36776283820.94
$ws.Range("U34").Value2 = 31619201138.78

$ws.Range("U36").Value2 = 567052855.599998
$ws.Range("U37").Value2 = 18841702314.3
$ws.Range("U38").Value2 = 9735751180.34
$ws.Range("U39").Value2 = 9587040690.81
$ws.Range("U40").Value2 = 5008188970.76
